$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "BiteAction" -> "ZombieAttackAction" (two occurrences in the first
#    paragraph, about Zombie attacks).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("BiteAction", $false, $false, $false, $false, $false, `
    $true, 1, $false, "ZombieAttackAction", 2)

# ---------------------------------------------------------------------------
# 2. Rewrite the "Rising from the Dead" paragraph.
#    Old:  ", we add a method toCorpse() which converts the Human to a
#          Corpse object is the human is dead. We check if the human is
#          dead by using the isConscious() that the Human class inherits
#          from the Actor class so that we abide by the "
#    New:  ", we created a ZombieCorpse class which inherits Item class.
#          We further changed the ZombieAttackAction class to be able to
#          knock out humans and create ZombieCorpse instead of normal
#          corpse item. The ZombieCorpse uses the tick() method from its
#          parent class Item so that we abide by the "
# ---------------------------------------------------------------------------
$oldPart1 = ", we add a method toCorpse() which converts the Human to a " + `
    "Corpse object is the human is dead. We check if the human is dead " + `
    "by using the isConscious() that the Human class inherits from the " + `
    "Actor class so that we abide by the "
$newPart1 = ", we created a ZombieCorpse class which inherits Item class. " + `
    "We further changed the ZombieAttackAction class to be able to knock " + `
    "out humans and create ZombieCorpse instead of normal corpse item. " + `
    "The ZombieCorpse uses the tick() method from its parent class Item " + `
    "so that we abide by the "
$d.Content.Find.Execute($oldPart1, $false, $false, $false, $false, $false, `
    $true, 1, $false, $newPart1, 2)

# Tail of the same paragraph: drop the old "method creates a Corpse
# object..." sentence, keep just " principle. ".
$oldPart2 = " principle. The method creates a Corpse object at the " + `
    "location using the addItem() method of the Location class."
$newPart2 = " principle. "
$d.Content.Find.Execute($oldPart2, $false, $false, $false, $false, $false, `
    $true, 1, $false, $newPart2, 2)

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark: it used to sit in the Farmers/Crop
#    paragraph (inside the "...tick() method of the item CLASS..."
#    sentence) and now belongs at the very end of the "Rising from the
#    Dead" paragraph, right after "... DRY principle. ".
#
#    Word keeps only one "_GoBack" bookmark at a time, so re-adding it at
#    the new spot automatically removes it from its old location.
#
#    A collapsed (zero-length) range that sits on the very last character
#    position of a paragraph (immediately before its paragraph mark) is
#    mishandled by this host when used directly with Bookmarks.Add, so we
#    work around it: insert a throw-away marker character after the
#    target point, anchor the bookmark just before that marker (no longer
#    the paragraph's last position), then delete the marker again.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("class Item so that we abide by the DRY principle. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
$anchor.InsertAfter("#")

$bookmarkSpot = $d.Range($anchor.Start, $anchor.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$marker = $d.Range($anchor.Start, $anchor.Start + 1)
$marker.Text = ""

# ---------------------------------------------------------------------------
# 4. Tidy up the Farmers/Crop paragraph: after the bookmark removal the
#    sentence is contiguous again, so merge it back into a single run by
#    re-applying the (unchanged) text through Find/Replace.
# ---------------------------------------------------------------------------
$cropSentence = " reduces by 1 until it reaches 0 by using the tick() " + `
    "method of the item class. The tick() method informs the Crop class " + `
    "of the passage of time."
$d.Content.Find.Execute($cropSentence, $false, $false, $false, $false, `
    $false, $true, 1, $false, $cropSentence, 2)
